$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.251.99"
$ws.Range("E2").Value = "  +2.49%  "
$ws.Range("D3").Value = "1.907.54"
$ws.Range("E3").Value = "  +2.26%  "
$dStyle4 = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = $dStyle4
$ws.Range("E4").Value = "  -0.12%  "
$dStyle5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.01"
$ws.Range("D5").Style = $dStyle5
$ws.Range("E5").Value = "  +1.23%  "
$dStyle6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = $dStyle6
$ws.Range("E6").Value = "  +0.01%  "
$dStyle7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4663"
$ws.Range("D7").Style = $dStyle7
$ws.Range("E7").Value = "  +1.13%  "
$dStyle8 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3964"
$ws.Range("D8").Style = $dStyle8
$ws.Range("E8").Value = "  +2.45%  "
$ws.Range("E9").Value = "  +1.10%  "
$dStyle10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07976"
$ws.Range("D10").Style = $dStyle10
$ws.Range("E10").Value = "  +1.53%  "
$dStyle11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.006"
$ws.Range("D11").Style = $dStyle11
$ws.Range("E11").Value = "  +3.23%  "
$dStyle12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.37"
$ws.Range("D12").Style = $dStyle12
$ws.Range("E12").Value = "  +2.15%  "
$ws.Range("D13").Value = "1.883.78"
$ws.Range("E13").Value = "  -0.39%  "
$dStyle14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.154"
$ws.Range("D14").Style = $dStyle14
$ws.Range("E14").Value = "  +2.62%  "
$dStyle15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.795"
$ws.Range("D15").Style = $dStyle15
$ws.Range("E15").Value = "  +1.90%  "
$ws.Range("E16").Value = "  +0.28%  "
$dStyle17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.94"
$ws.Range("D17").Style = $dStyle17
$ws.Range("E17").Value = "  +1.00%  "
$dStyle18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("D18").Style = $dStyle18
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("E19").Value = "  +1.32%  "
$dStyle20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.21"
$ws.Range("D20").Style = $dStyle20
$ws.Range("E20").Value = "  +2.52%  "
$dStyle21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").Style = $dStyle21
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "29.253.28"
$ws.Range("E22").Value = "  +2.43%  "
$dStyle23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.366"
$ws.Range("D23").Style = $dStyle23
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").Value = "2.138.70"
$ws.Range("E25").Value = "  +1.91%  "
$ws.Range("E26").Value = "  -2.44%  "
$dStyle27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.72"
$ws.Range("D27").Style = $dStyle27
$ws.Range("E27").Value = "  +2.91%  "
$dStyle28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.57"
$ws.Range("D28").Style = $dStyle28
$ws.Range("E28").Value = "  +1.75%  "
$dStyle29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.892"
$ws.Range("D29").Style = $dStyle29
$ws.Range("E29").Value = "  +1.97%  "
$dStyle30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.006"
$ws.Range("D30").Style = $dStyle30
$ws.Range("E30").Value = "  +1.14%  "
$dStyle31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.58"
$ws.Range("D31").Style = $dStyle31
$ws.Range("E31").Value = "  +0.45%  "
$dStyle32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09444"
$ws.Range("D32").Style = $dStyle32
$ws.Range("E32").Value = "  +1.46%  "
$dStyle33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9248"
$ws.Range("D33").Style = $dStyle33
$ws.Range("E33").Value = "  +0.85%  "
$dStyle34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.358"
$ws.Range("D34").Style = $dStyle34
$ws.Range("E34").Value = "  +1.79%  "
$dStyle35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.348"
$ws.Range("D35").Style = $dStyle35
$ws.Range("E35").Value = "  +1.25%  "
$dStyle36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.262"
$ws.Range("D36").Style = $dStyle36
$ws.Range("E36").Value = "  -1.87%  "
$dStyle37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05863"
$ws.Range("D37").Style = $dStyle37
$ws.Range("E37").Value = "  +1.51%  "
$dStyle38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.173"
$ws.Range("D38").Style = $dStyle38
$ws.Range("E38").Value = "  +2.05%  "
$dStyle39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02112"
$ws.Range("D39").Style = $dStyle39
$ws.Range("E39").Value = "  +2.13%  "
$dStyle40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.998"
$ws.Range("D40").Style = $dStyle40
$ws.Range("E40").Value = "  +3.60%  "
$dStyle41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5759"
$ws.Range("D41").Style = $dStyle41
$ws.Range("E41").Value = "  +2.49%  "
$ws.Range("E42").Value = "  +1.90%  "
$dStyle43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.04"
$ws.Range("D43").Style = $dStyle43
$ws.Range("E43").Value = "  +2.94%  "
$dStyle44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.03"
$ws.Range("D44").Style = $dStyle44
$ws.Range("E44").Value = "  +2.32%  "
$dStyle45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5440"
$ws.Range("D45").Style = $dStyle45
$ws.Range("E45").Value = "  +2.81%  "
$dStyle46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.229"
$ws.Range("D46").Style = $dStyle46
$ws.Range("E46").Value = "  +3.96%  "
$dStyle47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07091"
$ws.Range("D47").Style = $dStyle47
$ws.Range("E47").Value = "  -1.03%  "
$dStyle48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.889"
$ws.Range("D48").Style = $dStyle48
$ws.Range("E48").Value = "  +3.16%  "
$dStyle49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.580"
$ws.Range("D49").Style = $dStyle49
$ws.Range("E49").Value = "  +6.80%  "
$dStyle50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.09"
$ws.Range("D50").Style = $dStyle50
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("E51").Value = "  -5.26%  "
